# fix resource url mapping for geojson and update dataset_json_lookup
#
# The lookup table had a row mapping "url" -> "field_link_api" (with is_dataset
# = FALSE) that is no longer valid/needed. Remove that entire row (row 14),
# which shifts all subsequent rows up by one and drops the now-unused
# "field_link_api" shared string automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 currently holds: A="url", B="field_link_api", C=FALSE
$ws.Rows.Item(14).Delete()

# Reflect the new selection/active cell left behind after the edit.
$ws.Range("D27").Select()
